# Auto-generated edit script applying market-data value updates
# described by the commit diff across multiple sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 48.3
$ws.Range("I11").Value = 48.3
$ws.Range("K11").Value = 48.3
$ws.Range("M11").Value = 91.7
# Row 15
$ws.Range("H15").Value = 1101.8077
$ws.Range("I15").Value = 1101.8077
$ws.Range("K15").Value = 3305.4231
$ws.Range("M15").Value = -3136.4231
# Row 32
$ws.Range("H32").Value = 966.2
$ws.Range("I32").Value = 1500
$ws.Range("J32").Value = 772.0909
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 772.0909
$ws.Range("M32").Value = -1174
$ws.Range("N32").Value = -1424.0909
# Row 51
$ws.Range("H51").Value = 2328.6584
$ws.Range("I51").Value = 1527.2142
$ws.Range("J51").Value = 2744.2222
$ws.Range("K51").Value = 1527.2142
$ws.Range("L51").Value = 2744.2222
$ws.Range("M51").Value = -1043.2142
$ws.Range("N51").Value = -3712.2222
# Row 64
$ws.Range("H64").Value = 37040412
$ws.Range("I64").Value = 125002056
$ws.Range("J64").Value = 3929.4736
$ws.Range("K64").Value = 125002056
$ws.Range("L64").Value = 3929.4736
$ws.Range("M64").Value = -125001808
$ws.Range("N64").Value = -4425.473599999999
# Row 67
$ws.Range("H67").Value = 37040412
$ws.Range("I67").Value = 125002056
$ws.Range("J67").Value = 3929.4736
$ws.Range("K67").Value = 125002056
$ws.Range("L67").Value = 3929.4736
$ws.Range("M67").Value = -125001198
$ws.Range("N67").Value = -5645.473599999999
# Row 70
$ws.Range("H70").Value = 2940.4614
$ws.Range("I70").Value = 8067.3335
$ws.Range("J70").Value = 2271.739
$ws.Range("K70").Value = 24202.0005
$ws.Range("L70").Value = 6815.217000000001
$ws.Range("M70").Value = -23932.0005
$ws.Range("N70").Value = -7355.217000000001
# Row 73
$ws.Range("H73").Value = 2940.4614
$ws.Range("I73").Value = 8067.3335
$ws.Range("J73").Value = 2271.739
$ws.Range("K73").Value = 24202.0005
$ws.Range("L73").Value = 6815.217000000001
$ws.Range("M73").Value = -23266.0005
$ws.Range("N73").Value = -8687.217000000001
# Row 74
$ws.Range("H74").Value = 2988.25
$ws.Range("I74").Value = 2789
$ws.Range("J74").Value = 3187.5
$ws.Range("K74").Value = 2789
$ws.Range("L74").Value = 3187.5
$ws.Range("M74").Value = -1853
$ws.Range("N74").Value = -5059.5
# Row 77
$ws.Range("H77").Value = 2988.25
$ws.Range("I77").Value = 2789
$ws.Range("J77").Value = 3187.5
$ws.Range("K77").Value = 13945
$ws.Range("L77").Value = 15937.5
$ws.Range("M77").Value = -9265
$ws.Range("N77").Value = -25297.5
# Row 100
$ws.Range("H100").Value = 33161.938
$ws.Range("I100").Value = 51363.5
$ws.Range("K100").Value = 51363.5
$ws.Range("M100").Value = -50822.5

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 2203.9033
$ws.Range("I97").Value = 600.4167
$ws.Range("J97").Value = 7701.5713
$ws.Range("K97").Value = 600.4167
$ws.Range("L97").Value = 7701.5713
$ws.Range("M97").Value = -104.4167
$ws.Range("N97").Value = -8693.5713
# Row 110
$ws.Range("H110").Value = 6072.364
$ws.Range("I110").Value = 4678.3
$ws.Range("J110").Value = 20013
$ws.Range("K110").Value = 4678.3
$ws.Range("L110").Value = 20013
$ws.Range("M110").Value = -2633.3
$ws.Range("N110").Value = -24103

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2794.9473
$ws.Range("I86").Value = 1738.5555
$ws.Range("J86").Value = 3745.7
$ws.Range("K86").Value = 1738.5555
$ws.Range("L86").Value = 3745.7
$ws.Range("M86").Value = -615.5554999999999
$ws.Range("N86").Value = -5991.7
# Row 89
$ws.Range("H89").Value = 2794.9473
$ws.Range("I89").Value = 1738.5555
$ws.Range("J89").Value = 3745.7
$ws.Range("K89").Value = 8692.7775
$ws.Range("L89").Value = 18728.5
$ws.Range("M89").Value = -3076.7775
$ws.Range("N89").Value = -29960.5
# Row 94
$ws.Range("H94").Value = 1230.075
$ws.Range("I94").Value = 799.3929000000001
$ws.Range("J94").Value = 2235
$ws.Range("K94").Value = 799.3929000000001
$ws.Range("L94").Value = 2235
$ws.Range("M94").Value = -348.3929000000001
$ws.Range("N94").Value = -3137
# Row 99
$ws.Range("H99").Value = 2402.2
$ws.Range("J99").Value = 3337
$ws.Range("L99").Value = 3337
$ws.Range("N99").Value = -6333
# Row 107
$ws.Range("H107").Value = 1601.4166
$ws.Range("I107").Value = 901.375
$ws.Range("K107").Value = 901.375
$ws.Range("M107").Value = 1018.625

$ws = $wb.Worksheets.Item("CRP")
# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# Row 132
$ws.Range("H132").Value = 26317744
$ws.Range("I132").Value = 33334810
$ws.Range("J132").Value = 3747.25
$ws.Range("K132").Value = 100004430
$ws.Range("L132").Value = 11241.75
$ws.Range("M132").Value = -100001900
$ws.Range("N132").Value = -16301.75
# Row 134
$ws.Range("H134").Value = 2752.3635
$ws.Range("I134").Value = 2830.2856
$ws.Range("J134").Value = 2616
$ws.Range("K134").Value = 8490.856800000001
$ws.Range("L134").Value = 7848
$ws.Range("M134").Value = -5955.856800000001
$ws.Range("N134").Value = -12918

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 545.8333
$ws.Range("J5").Value = 900
$ws.Range("L5").Value = 2700
$ws.Range("N5").Value = -2924
# Row 49
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 2600
$ws.Range("K49").Value = 3000
$ws.Range("L49").Value = 7800
$ws.Range("M49").Value = -2844
$ws.Range("N49").Value = -8112
# Row 92
$ws.Range("H92").Value = 2451.5
$ws.Range("I92").Value = 3900
$ws.Range("J92").Value = 1003
$ws.Range("K92").Value = 11700
$ws.Range("L92").Value = 3009
$ws.Range("M92").Value = -10452
$ws.Range("N92").Value = -5505
# Row 135
$ws.Range("H135").Value = 545.8333
$ws.Range("J135").Value = 900
$ws.Range("L135").Value = 8100
$ws.Range("N135").Value = -13170

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1121.7858
$ws.Range("I22").Value = 1230.1666
$ws.Range("J22").Value = 1040.5
$ws.Range("K22").Value = 1230.1666
$ws.Range("L22").Value = 1040.5
$ws.Range("M22").Value = -935.1666
$ws.Range("N22").Value = -1630.5
# Row 27
$ws.Range("H27").Value = 1121.7858
$ws.Range("I27").Value = 1230.1666
$ws.Range("J27").Value = 1040.5
$ws.Range("K27").Value = 1230.1666
$ws.Range("L27").Value = 1040.5
$ws.Range("M27").Value = -1123.1666
$ws.Range("N27").Value = -1254.5
# Row 46
$ws.Range("H46").Value = 3634.125
$ws.Range("I46").Value = 5735.5
$ws.Range("J46").Value = 2933.6667
$ws.Range("K46").Value = 5735.5
$ws.Range("L46").Value = 2933.6667
$ws.Range("M46").Value = -5547.5
$ws.Range("N46").Value = -3309.6667
# Row 64
$ws.Range("H64").Value = 34142.855
$ws.Range("J64").Value = 34142.855
$ws.Range("L64").Value = 34142.855
$ws.Range("N64").Value = -34592.855
# Row 67
$ws.Range("H67").Value = 34142.855
$ws.Range("J67").Value = 34142.855
$ws.Range("L67").Value = 34142.855
$ws.Range("N67").Value = -35702.855
# Row 82
$ws.Range("H82").Value = 1500
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1500
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1500
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2222
# Row 85
$ws.Range("H85").Value = 1500
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1500
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1500
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3996
# Row 93
$ws.Range("H93").Value = 1610.2142
$ws.Range("I93").Value = 1500.5
$ws.Range("J93").Value = 1692.5
$ws.Range("K93").Value = 1500.5
$ws.Range("L93").Value = 1692.5
$ws.Range("M93").Value = -252.5
$ws.Range("N93").Value = -4188.5

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2720.1304
$ws.Range("I62").Value = 2751.3333
$ws.Range("J62").Value = 2661.625
$ws.Range("K62").Value = 2751.3333
$ws.Range("L62").Value = 2661.625
$ws.Range("M62").Value = -2127.3333
$ws.Range("N62").Value = -3909.625
# Row 63
$ws.Range("H63").Value = 7000
$ws.Range("J63").Value = 7000
$ws.Range("L63").Value = 7000
$ws.Range("N63").Value = -8248
# Row 65
$ws.Range("H65").Value = 2720.1304
$ws.Range("I65").Value = 2751.3333
$ws.Range("J65").Value = 2661.625
$ws.Range("K65").Value = 13756.6665
$ws.Range("L65").Value = 13308.125
$ws.Range("M65").Value = -10636.6665
$ws.Range("N65").Value = -19548.125
# Row 66
$ws.Range("H66").Value = 7000
$ws.Range("J66").Value = 7000
$ws.Range("L66").Value = 21000
$ws.Range("N66").Value = -27240
# Row 81
$ws.Range("H81").Value = 1400.5
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 1700.6666
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 3401.3332
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -5523.3332
# Row 84
$ws.Range("H84").Value = 1400.5
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 1700.6666
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 17006.666
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -27614.666
